# Update cryptos list with latest prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage each changed cell as a text-literal formula so Excel keeps the value as
# text (these price/volume columns are plain text in the workbook, and some of
# the new values - e.g. "10.50", "0.999" - look numeric and would otherwise be
# auto-coerced to a Number by a plain .Value assignment).
$ws.Range("D2").Formula = '="64.392.19"'
$ws.Range("E2").Formula = '="  +0.43%  "'
$ws.Range("D3").Formula = '="3.161.65"'
$ws.Range("E3").Formula = '="  -0.16%  "'
$ws.Range("E4").Formula = '="  +0.04%  "'
$ws.Range("D5").Formula = '="571.65"'
$ws.Range("E5").Formula = '="  +0.50%  "'
$ws.Range("D6").Formula = '="164.01"'
$ws.Range("E6").Formula = '="  -2.76%  "'
$ws.Range("E7").Formula = '="  +0.09%  "'
$ws.Range("D8").Formula = '="0.578"'
$ws.Range("E8").Formula = '="  -4.66%  "'
$ws.Range("E9").Formula = '="  -2.49%  "'
$ws.Range("D10").Formula = '="6.62"'
$ws.Range("E10").Formula = '="  -1.14%  "'
$ws.Range("D11").Formula = '="0.383"'
$ws.Range("E11").Formula = '="  -0.11%  "'
$ws.Range("D12").Formula = '="3.712.69"'
$ws.Range("E12").Formula = '="  -0.10%  "'
$ws.Range("E13").Formula = '="  -0.86%  "'
$ws.Range("D14").Formula = '="64.429.16"'
$ws.Range("E14").Formula = '="  +0.40%  "'
$ws.Range("D15").Formula = '="25.21"'
$ws.Range("E15").Formula = '="  -0.29%  "'
$ws.Range("D16").Formula = '="3.160.94"'
$ws.Range("E16").Formula = '="  +0.23%  "'
$ws.Range("E17").Formula = '="  -2.11%  "'
$ws.Range("D18").Formula = '="405.10"'
$ws.Range("E18").Formula = '="  -2.50%  "'
$ws.Range("D19").Formula = '="12.72"'
$ws.Range("E19").Formula = '="  -0.67%  "'
$ws.Range("D20").Formula = '="5.24"'
$ws.Range("E20").Formula = '="  -1.76%  "'
$ws.Range("D21").Formula = '="7.12"'
$ws.Range("E21").Formula = '="  +0.59%  "'
$ws.Range("E22").Formula = '="  +0.21%  "'
$ws.Range("D23").Formula = '="68.68"'
$ws.Range("E23").Formula = '="  -1.89%  "'
$ws.Range("E24").Formula = '="  -0.98%  "'
$ws.Range("D25").Formula = '="0.192"'
$ws.Range("E25").Formula = '="  -4.64%  "'
$ws.Range("E26").Formula = '="  -3.83%  "'
$ws.Range("D27").Formula = '="8.81"'
$ws.Range("E27").Formula = '="  +0.72%  "'
$ws.Range("D28").Formula = '="0.998"'
$ws.Range("E28").Formula = '="  +0.18%  "'
$ws.Range("D29").Formula = '="1.81"'
$ws.Range("E29").Formula = '="  -0.66%  "'
$ws.Range("D30").Formula = '="21.19"'
$ws.Range("E30").Formula = '="  -2.55%  "'
$ws.Range("D31").Formula = '="6.31"'
$ws.Range("E31").Formula = '="  -0.21%  "'
$ws.Range("D32").Formula = '="4.85"'
$ws.Range("E32").Formula = '="  -3.27%  "'
$ws.Range("B33").Formula = '="Monero"'
$ws.Range("C33").Formula = '="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"'
$ws.Range("D33").Formula = '="156.76"'
$ws.Range("E33").Formula = '="  +1.18%  "'
$ws.Range("B34").Formula = '="Fetch.AI"'
$ws.Range("C34").Formula = '="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"'
$ws.Range("D34").Formula = '="1.12"'
$ws.Range("E34").Formula = '="  -0.98%  "'
$ws.Range("E35").Formula = '="  -2.37%  "'
$ws.Range("B36").Formula = '="Maker"'
$ws.Range("C36").Formula = '="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"'
$ws.Range("D36").Formula = '="2.674.72"'
$ws.Range("E36").Formula = '="  -0.88%  "'
$ws.Range("B37").Formula = '="Stacks"'
$ws.Range("C37").Formula = '="https://coinranking.com/coin/mMPrMcB7+stacks-stx"'
$ws.Range("D37").Formula = '="1.68"'
$ws.Range("E37").Formula = '="  -0.75%  "'
$ws.Range("D38").Formula = '="23.83"'
$ws.Range("E38").Formula = '="  -2.76%  "'
$ws.Range("D39").Formula = '="4.09"'
$ws.Range("E39").Formula = '="  -1.97%  "'
$ws.Range("D40").Formula = '="0.695"'
$ws.Range("E40").Formula = '="  -1.87%  "'
$ws.Range("D41").Formula = '="0.0617"'
$ws.Range("E41").Formula = '="  -1.09%  "'
$ws.Range("E42").Formula = '="  -2.94%  "'
$ws.Range("B43").Formula = '="Bittensor"'
$ws.Range("C43").Formula = '="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"'
$ws.Range("D43").Formula = '="289.40"'
$ws.Range("E43").Formula = '="  -1.82%  "'
$ws.Range("B44").Formula = '="VeChain"'
$ws.Range("C44").Formula = '="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"'
$ws.Range("D44").Formula = '="0.0256"'
$ws.Range("E44").Formula = '="  -2.12%  "'
$ws.Range("B45").Formula = '="InjectiveProtocol"'
$ws.Range("C45").Formula = '="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"'
$ws.Range("D45").Formula = '="21.29"'
$ws.Range("E45").Formula = '="  -2.29%  "'
$ws.Range("D46").Formula = '="0.999"'
$ws.Range("E46").Formula = '="  +0.00%  "'
$ws.Range("D47").Formula = '="0.0981"'
$ws.Range("E47").Formula = '="  -0.96%  "'
$ws.Range("D48").Formula = '="10.50"'
$ws.Range("E48").Formula = '="  +0.70%  "'
$ws.Range("D49").Formula = '="1.90"'
$ws.Range("E49").Formula = '="  -7.17%  "'
$ws.Range("D50").Formula = '="5.70"'
$ws.Range("E50").Formula = '="  -1.33%  "'
$ws.Range("D51").Formula = '="0.876"'
$ws.Range("E51").Formula = '="  -6.70%  "'

# Convert the staged formulas back into static text values in one shot, without
# touching any cell styles/number formats.
$ws.Range("B2:E51").Copy()
$ws.Range("B2:E51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
